$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 to hold what was previously row 3's data
# (Officer Registration ID 1's record now reflects the successful
# registration details, and the old row 2 registration record is removed).
$ws.Range("B2").Value = "T1234567J"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Successful"
$ws.Range("E2").Value = 45769.794563912037

# Remove the now-duplicate row (old row 3)
$ws.Rows.Item(3).Delete()

# Reflect the final cell selection saved with the workbook
[void]$ws.Range("B8").Select()
